# Best_XI.xlsx update: rename sheets (prefix "g"), switch the active tab
# from the 1st sheet ("26"/"g26") to the 2nd sheet ("27"/"g27"), and move
# the remembered selection on sheets 2 and 3.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Rename the three sheets, keeping their relative order: "26"->"g26",
# "27"->"g27", "28"->"g28".
$ws1.Name = "g26"
$ws2.Name = "g27"
$ws3.Name = "g28"

# Sheet 1 ("g26") keeps its old selection (F12); it just stops being the
# active tab, which happens automatically once another sheet is activated.
[void]$ws1.Range("F12").Select()

# Sheet 3 ("g28") remembers a new selection (E28) but is not the active tab.
[void]$ws3.Range("E28").Select()

# Sheet 2 ("g27") becomes the active tab, with a new remembered selection
# (E11). Doing this last makes it the active sheet in the saved workbook.
[void]$ws2.Activate()
[void]$ws2.Range("E11").Select()
